$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-26 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-27 Saturday", 2) | Out-Null
$d.Content.Find.Execute("37+31=68", $true, $false, $false, $false, $false, $true, 1, $false, "84-30=54", 2) | Out-Null
$d.Content.Find.Execute("58-22=36", $true, $false, $false, $false, $false, $true, 1, $false, "48+36=84", 2) | Out-Null
$d.Content.Find.Execute("79-54=25", $true, $false, $false, $false, $false, $true, 1, $false, "63-54=9", 2) | Out-Null
$d.Content.Find.Execute("33+31=64", $true, $false, $false, $false, $false, $true, 1, $false, "76-58=18", 2) | Out-Null
$d.Content.Find.Execute("85-38=47", $true, $false, $false, $false, $false, $true, 1, $false, "72-25=47", 2) | Out-Null
$d.Content.Find.Execute("52+13=65", $true, $false, $false, $false, $false, $true, 1, $false, "95-73=22", 2) | Out-Null
$d.Content.Find.Execute("19+2=21", $true, $false, $false, $false, $false, $true, 1, $false, "46+49=95", 2) | Out-Null
$d.Content.Find.Execute("71-16=55", $true, $false, $false, $false, $false, $true, 1, $false, "9+17=26", 2) | Out-Null
$d.Content.Find.Execute("60-50=10", $true, $false, $false, $false, $false, $true, 1, $false, "81+4=85", 2) | Out-Null
$d.Content.Find.Execute("81-0=81", $true, $false, $false, $false, $false, $true, 1, $false, "14-0=14", 2) | Out-Null
$d.Content.Find.Execute("58+33=91", $true, $false, $false, $false, $false, $true, 1, $false, "35+7=42", 2) | Out-Null
$d.Content.Find.Execute("87-12=75", $true, $false, $false, $false, $false, $true, 1, $false, "58-4=54", 2) | Out-Null
$d.Content.Find.Execute("28+42=70", $true, $false, $false, $false, $false, $true, 1, $false, "36+19=55", 2) | Out-Null
$d.Content.Find.Execute("30+63=93", $true, $false, $false, $false, $false, $true, 1, $false, "23+72=95", 2) | Out-Null
$d.Content.Find.Execute("51-13=38", $true, $false, $false, $false, $false, $true, 1, $false, "85-49=36", 2) | Out-Null
$d.Content.Find.Execute("67-61=6", $true, $false, $false, $false, $false, $true, 1, $false, "22-13=9", 2) | Out-Null
$d.Content.Find.Execute("63-26=37", $true, $false, $false, $false, $false, $true, 1, $false, "1+12=13", 2) | Out-Null
$d.Content.Find.Execute("90-15=75", $true, $false, $false, $false, $false, $true, 1, $false, "5+16=21", 2) | Out-Null
$d.Content.Find.Execute("77-16=61", $true, $false, $false, $false, $false, $true, 1, $false, "38+25=63", 2) | Out-Null
$d.Content.Find.Execute("60+25=85", $true, $false, $false, $false, $false, $true, 1, $false, "58-0=58", 2) | Out-Null
$d.Content.Find.Execute("28-22=6", $true, $false, $false, $false, $false, $true, 1, $false, "6+93=99", 2) | Out-Null
$d.Content.Find.Execute("78+21=99", $true, $false, $false, $false, $false, $true, 1, $false, "95-62=33", 2) | Out-Null
$d.Content.Find.Execute("56-4=52", $true, $false, $false, $false, $false, $true, 1, $false, "73+11=84", 2) | Out-Null
$d.Content.Find.Execute("97-38=59", $true, $false, $false, $false, $false, $true, 1, $false, "88-63=25", 2) | Out-Null
$d.Content.Find.Execute("8+86=94", $true, $false, $false, $false, $false, $true, 1, $false, "89-10=79", 2) | Out-Null
$d.Content.Find.Execute("76+17=93", $true, $false, $false, $false, $false, $true, 1, $false, "64+14=78", 2) | Out-Null
$d.Content.Find.Execute("31+45=76", $true, $false, $false, $false, $false, $true, 1, $false, "32+40=72", 2) | Out-Null
$d.Content.Find.Execute("53-20=33", $true, $false, $false, $false, $false, $true, 1, $false, "29+54=83", 2) | Out-Null
$d.Content.Find.Execute("7+51=58", $true, $false, $false, $false, $false, $true, 1, $false, "18+27=45", 2) | Out-Null
$d.Content.Find.Execute("50+47=97", $true, $false, $false, $false, $false, $true, 1, $false, "76+22=98", 2) | Out-Null
$d.Content.Find.Execute("18+35=53", $true, $false, $false, $false, $false, $true, 1, $false, "7+8=15", 2) | Out-Null
$d.Content.Find.Execute("48-31=17", $true, $false, $false, $false, $false, $true, 1, $false, "79-47=32", 2) | Out-Null
$d.Content.Find.Execute("82-53=29", $true, $false, $false, $false, $false, $true, 1, $false, "24+41=65", 2) | Out-Null
$d.Content.Find.Execute("53-53=0", $true, $false, $false, $false, $false, $true, 1, $false, "32-32=0", 2) | Out-Null
$d.Content.Find.Execute("91-17=74", $true, $false, $false, $false, $false, $true, 1, $false, "4+60=64", 2) | Out-Null
$d.Content.Find.Execute("30+24=54", $true, $false, $false, $false, $false, $true, 1, $false, "21+62=83", 2) | Out-Null
$d.Content.Find.Execute("65-64=1", $true, $false, $false, $false, $false, $true, 1, $false, "10+64=74", 2) | Out-Null
$d.Content.Find.Execute("9+70=79", $true, $false, $false, $false, $false, $true, 1, $false, "19+72=91", 2) | Out-Null
$d.Content.Find.Execute("30-29=1", $true, $false, $false, $false, $false, $true, 1, $false, "63-27=36", 2) | Out-Null
$d.Content.Find.Execute("12+27=39", $true, $false, $false, $false, $false, $true, 1, $false, "42+48=90", 2) | Out-Null
$d.Content.Find.Execute("31+19=50", $true, $false, $false, $false, $false, $true, 1, $false, "15+78=93", 2) | Out-Null
$d.Content.Find.Execute("24-19=5", $true, $false, $false, $false, $false, $true, 1, $false, "58-15=43", 2) | Out-Null
$d.Content.Find.Execute("18+26=44", $true, $false, $false, $false, $false, $true, 1, $false, "13+43=56", 2) | Out-Null
$d.Content.Find.Execute("18+0=18", $true, $false, $false, $false, $false, $true, 1, $false, "93-80=13", 2) | Out-Null
$d.Content.Find.Execute("9+2=11", $true, $false, $false, $false, $false, $true, 1, $false, "81+9=90", 2) | Out-Null
$d.Content.Find.Execute("21+9=30", $true, $false, $false, $false, $false, $true, 1, $false, "76+13=89", 2) | Out-Null
$d.Content.Find.Execute("93-8=85", $true, $false, $false, $false, $false, $true, 1, $false, "61+6=67", 2) | Out-Null
$d.Content.Find.Execute("38+51=89", $true, $false, $false, $false, $false, $true, 1, $false, "24-12=12", 2) | Out-Null
$d.Content.Find.Execute("26-24=2", $true, $false, $false, $false, $false, $true, 1, $false, "26-20=6", 2) | Out-Null
$d.Content.Find.Execute("2+87=89", $true, $false, $false, $false, $false, $true, 1, $false, "79+8=87", 2) | Out-Null
$d.Content.Find.Execute("5+14=19", $true, $false, $false, $false, $false, $true, 1, $false, "2+2=4", 2) | Out-Null
$d.Content.Find.Execute("2+20=22", $true, $false, $false, $false, $false, $true, 1, $false, "35+30=65", 2) | Out-Null
$d.Content.Find.Execute("93-62=31", $true, $false, $false, $false, $false, $true, 1, $false, "39+8=47", 2) | Out-Null
$d.Content.Find.Execute("74+13=87", $true, $false, $false, $false, $false, $true, 1, $false, "36+18=54", 2) | Out-Null
$d.Content.Find.Execute("0+45=45", $true, $false, $false, $false, $false, $true, 1, $false, "70+4=74", 2) | Out-Null
$d.Content.Find.Execute("57+34=91", $true, $false, $false, $false, $false, $true, 1, $false, "7+39=46", 2) | Out-Null
$d.Content.Find.Execute("46-8=38", $true, $false, $false, $false, $false, $true, 1, $false, "20-9=11", 2) | Out-Null
$d.Content.Find.Execute("75-50=25", $true, $false, $false, $false, $false, $true, 1, $false, "71-65=6", 2) | Out-Null
$d.Content.Find.Execute("48-2=46", $true, $false, $false, $false, $false, $true, 1, $false, "55-35=20", 2) | Out-Null
$d.Content.Find.Execute("15+46=61", $true, $false, $false, $false, $false, $true, 1, $false, "84-50=34", 2) | Out-Null
$d.Content.Find.Execute("86-53=33", $true, $false, $false, $false, $false, $true, 1, $false, "31-31=0", 2) | Out-Null
$d.Content.Find.Execute("56+6=62", $true, $false, $false, $false, $false, $true, 1, $false, "64-10=54", 2) | Out-Null
$d.Content.Find.Execute("62-16=46", $true, $false, $false, $false, $false, $true, 1, $false, "80-28=52", 2) | Out-Null
$d.Content.Find.Execute("58-13=45", $true, $false, $false, $false, $false, $true, 1, $false, "76+16=92", 2) | Out-Null
$d.Content.Find.Execute("50-16=34", $true, $false, $false, $false, $false, $true, 1, $false, "35+32=67", 2) | Out-Null
$d.Content.Find.Execute("39+34=73", $true, $false, $false, $false, $false, $true, 1, $false, "80-51=29", 2) | Out-Null
$d.Content.Find.Execute("13-3=10", $true, $false, $false, $false, $false, $true, 1, $false, "47+52=99", 2) | Out-Null
$d.Content.Find.Execute("14-6=8", $true, $false, $false, $false, $false, $true, 1, $false, "44+36=80", 2) | Out-Null
$d.Content.Find.Execute("36+50=86", $true, $false, $false, $false, $false, $true, 1, $false, "99-54=45", 2) | Out-Null
$d.Content.Find.Execute("42+9=51", $true, $false, $false, $false, $false, $true, 1, $false, "18+46=64", 2) | Out-Null
$d.Content.Find.Execute("43+0=43", $true, $false, $false, $false, $false, $true, 1, $false, "39-33=6", 2) | Out-Null
$d.Content.Find.Execute("26+5=31", $true, $false, $false, $false, $false, $true, 1, $false, "9+22=31", 2) | Out-Null
$d.Content.Find.Execute("64+26=90", $true, $false, $false, $false, $false, $true, 1, $false, "38+23=61", 2) | Out-Null
$d.Content.Find.Execute("74+8=82", $true, $false, $false, $false, $false, $true, 1, $false, "63+13=76", 2) | Out-Null
$d.Content.Find.Execute("17+15=32", $true, $false, $false, $false, $false, $true, 1, $false, "44-3=41", 2) | Out-Null
$d.Content.Find.Execute("81-27=54", $true, $false, $false, $false, $false, $true, 1, $false, "34-3=31", 2) | Out-Null
$d.Content.Find.Execute("0+62=62", $true, $false, $false, $false, $false, $true, 1, $false, "14+70=84", 2) | Out-Null
$d.Content.Find.Execute("43+23=66", $true, $false, $false, $false, $false, $true, 1, $false, "8+81=89", 2) | Out-Null
$d.Content.Find.Execute("97-94=3", $true, $false, $false, $false, $false, $true, 1, $false, "41+53=94", 2) | Out-Null
$d.Content.Find.Execute("88-35=53", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=7", 2) | Out-Null
$d.Content.Find.Execute("94-12=82", $true, $false, $false, $false, $false, $true, 1, $false, "4+14=18", 2) | Out-Null
$d.Content.Find.Execute("22+47=69", $true, $false, $false, $false, $false, $true, 1, $false, "15+49=64", 2) | Out-Null
$d.Content.Find.Execute("7+14=21", $true, $false, $false, $false, $false, $true, 1, $false, "81-60=21", 2) | Out-Null
$d.Content.Find.Execute("91-56=35", $true, $false, $false, $false, $false, $true, 1, $false, "19-17=2", 2) | Out-Null
$d.Content.Find.Execute("78-3=75", $true, $false, $false, $false, $false, $true, 1, $false, "61+24=85", 2) | Out-Null
$d.Content.Find.Execute("87-45=42", $true, $false, $false, $false, $false, $true, 1, $false, "46+28=74", 2) | Out-Null
$d.Content.Find.Execute("33+24=57", $true, $false, $false, $false, $false, $true, 1, $false, "99-94=5", 2) | Out-Null
$d.Content.Find.Execute("28+48=76", $true, $false, $false, $false, $false, $true, 1, $false, "22+61=83", 2) | Out-Null
$d.Content.Find.Execute("7-0=7", $true, $false, $false, $false, $false, $true, 1, $false, "78-49=29", 2) | Out-Null
$d.Content.Find.Execute("58+32=90", $true, $false, $false, $false, $false, $true, 1, $false, "14+5=19", 2) | Out-Null
$d.Content.Find.Execute("79-9=70", $true, $false, $false, $false, $false, $true, 1, $false, "99-57=42", 2) | Out-Null
$d.Content.Find.Execute("3+73=76", $true, $false, $false, $false, $false, $true, 1, $false, "22+65=87", 2) | Out-Null
$d.Content.Find.Execute("51+4=55", $true, $false, $false, $false, $false, $true, 1, $false, "34-31=3", 2) | Out-Null
$d.Content.Find.Execute("14+34=48", $true, $false, $false, $false, $false, $true, 1, $false, "33-4=29", 2) | Out-Null
$d.Content.Find.Execute("94-68=26", $true, $false, $false, $false, $false, $true, 1, $false, "25+26=51", 2) | Out-Null
$d.Content.Find.Execute("60-34=26", $true, $false, $false, $false, $false, $true, 1, $false, "94-45=49", 2) | Out-Null
$d.Content.Find.Execute("58-56=2", $true, $false, $false, $false, $false, $true, 1, $false, "21+62=83", 2) | Out-Null
$d.Content.Find.Execute("34+7=41", $true, $false, $false, $false, $false, $true, 1, $false, "41-24=17", 2) | Out-Null
$d.Content.Find.Execute("79-13=66", $true, $false, $false, $false, $false, $true, 1, $false, "51-35=16", 2) | Out-Null
$d.Content.Find.Execute("16+30=46", $true, $false, $false, $false, $false, $true, 1, $false, "41+30=71", 2) | Out-Null
